$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")
for ($r=17; $r -le 42; $r++) {
    $b = $ws.Range("B$r").Text
    $c = $ws.Range("C$r").Text
    $d = $ws.Range("D$r").Text
    $e = $ws.Range("E$r").Text
    $f = $ws.Range("F$r").Text
    $h = $ws.Range("H$r").Text
    Write-Host "Row $r : B=[$b] C=[$c] D=[$d] E=[$e] F=[$f] H=[$h]"
}
